# Scene.xlsx — add two plain (non XML-mapped) columns "CanClone" and
# "ActorID" to the end of the existing XML table, fill their values for
# the 3 data rows, widen the new "CanClone" column, move the selection,
# and turn on an explicit (portrait / A4) page setup — matching the
# authored commit "modified configuration of npc up the bin file of server".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table by two columns (A1:I4 -> A1:K4). ListColumns.Add()
# appends one column to the right of the table and expands ref/autoFilter
# accordingly.
$colCanClone = $lo.ListColumns.Add()
$colActorID  = $lo.ListColumns.Add()

# Header labels. "ActorID" (column K) is written before "CanClone"
# (column J) so the shared-string table gets the two new entries in
# that order, matching the source workbook.
$ws.Range("K1").Value2 = "ActorID"
$ws.Range("J1").Value2 = "CanClone"

# Per-row values for the two new columns.
$ws.Range("J2").Value2 = 1
$ws.Range("K2").Value2 = 0

$ws.Range("J3").Value2 = 0
$ws.Range("K3").Value2 = 0

$ws.Range("J4").Value2 = 0
$ws.Range("K4").Value2 = 0

# New "CanClone" column (10th sheet column / J) gets an explicit width,
# matching the width already used by the other "14"-wide columns.
$ws.Columns.Item(10).ColumnWidth = $ws.Columns.Item(8).ColumnWidth

# Active selection moves to K9.
$null = $ws.Range("K9").Select()

# Explicit page setup (portrait, paper size 9 = A4).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
